# Edit slide 5 ("Solución Mediante ROI Frente - Boca"):
#  - move the first chart picture ("Imagen 3") further down the slide
#  - shrink the first value textbox and relabel it "BPM" -> "SpO2"
#  - move/resize the second picture ("Picture 2") up into the vacated gap
#  - add a new "BPM" textbox below the relocated picture (a duplicate of
#    the original label box, so it keeps the same run/body formatting)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)

$imagen3  = $s.Shapes.Item(2)   # "Imagen 3" picture
$cuadro5  = $s.Shapes.Item(3)   # "CuadroTexto 5" textbox (currently "BPM")
$picture2 = $s.Shapes.Item(4)   # "Picture 2" picture

# 4) Duplicate the "BPM" textbox first (while it still has its original
#    size/text) so the new shape inherits identical formatting, then move
#    it into place under the relocated "Picture 2".
$dup = $cuadro5.Duplicate().Item(1)
$dup.Name = "CuadroTexto 2"
$dup.Left = 31.92976377952756          # EMU 405508
$dup.Top = 240.91874015748033          # EMU 3059668
$dup.Width = 62.43913465826772         # EMU 792977
$dup.Height = 29.081259842519685       # EMU 369332
$dup.TextFrame.TextRange.Text = "BPM"

# 1) "Imagen 3" picture: slide it down (only Top changes)
$imagen3.Top = 279.5672454944882       # EMU 3550504 (was 1381994)

# 2) "CuadroTexto 5" textbox: narrower box, new caption
$cuadro5.Width = 62.43913465826772     # EMU 792977 (was 3154438)
$cuadro5.TextFrame.TextRange.Text = "SpO2"

# 3) "Picture 2": reposition/resize into the space vacated above
$picture2.Left = 31.92976377952756     # EMU 405508
$picture2.Top = 108.8184251968504      # EMU 1381994
$picture2.Width = 896.1404724409449    # EMU 11380984
$picture2.Height = 123.07944881889763  # EMU 1563109
